# Adds rows 128-181 to Sheet1, extending the collateral history table
# with a new valuation date block (two new portfolio/counterparty summaries,
# duplicated twice with the counterparty long-name variant), per the commit
# "gestion de collat modifs".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Same datetime number format used by the existing "Date" column (style applied
# to A2:A127) so new date cells share the existing style instead of creating a new one.
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Row 128
$ws.Cells.Item(128, 1).Value = 45916
$ws.Cells.Item(128, 1).NumberFormat = $dateFormat
$ws.Cells.Item(128, 2).Value = 900200
$ws.Cells.Item(128, 3).Value = 'BNPP'
$ws.Cells.Item(128, 4).Value = -44748309.06889844
$ws.Cells.Item(128, 5).Value = 0
$ws.Cells.Item(128, 6).Value = 0
$ws.Cells.Item(128, 7).Value = 'cash insuffisant (22 757 791.25)'

# Row 129
$ws.Cells.Item(129, 1).Value = 45916
$ws.Cells.Item(129, 1).NumberFormat = $dateFormat
$ws.Cells.Item(129, 2).Value = 900200
$ws.Cells.Item(129, 3).Value = 'CA'
$ws.Cells.Item(129, 4).Value = -16145545.642625
$ws.Cells.Item(129, 5).Value = 0
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 'cash insuffisant (16 145 545.64)'

# Row 130
$ws.Cells.Item(130, 1).Value = 45916
$ws.Cells.Item(130, 1).NumberFormat = $dateFormat
$ws.Cells.Item(130, 2).Value = 900200
$ws.Cells.Item(130, 3).Value = 'CEP'
$ws.Cells.Item(130, 4).Value = 110691.9541025327
$ws.Cells.Item(130, 5).Value = 0
$ws.Cells.Item(130, 6).Value = 0

# Row 131
$ws.Cells.Item(131, 1).Value = 45916
$ws.Cells.Item(131, 1).NumberFormat = $dateFormat
$ws.Cells.Item(131, 2).Value = 900200
$ws.Cells.Item(131, 3).Value = 'GIPB'
$ws.Cells.Item(131, 4).Value = -262844.824426329
$ws.Cells.Item(131, 5).Value = 0
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 'cash insuffisant (262 844.82)'

# Row 132
$ws.Cells.Item(132, 1).Value = 45916
$ws.Cells.Item(132, 1).NumberFormat = $dateFormat
$ws.Cells.Item(132, 2).Value = 900200
$ws.Cells.Item(132, 3).Value = 'GSOH'
$ws.Cells.Item(132, 4).Value = -1552465.916127797
$ws.Cells.Item(132, 5).Value = 0
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 'cash insuffisant (1 552 465.92)'

# Row 133
$ws.Cells.Item(133, 1).Value = 45916
$ws.Cells.Item(133, 1).NumberFormat = $dateFormat
$ws.Cells.Item(133, 2).Value = 900200
$ws.Cells.Item(133, 3).Value = 'JPMSE'
$ws.Cells.Item(133, 4).Value = -6725967.70460027
$ws.Cells.Item(133, 5).Value = 0
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 'cash insuffisant (6 725 967.70)'

# Row 134
$ws.Cells.Item(134, 1).Value = 45916
$ws.Cells.Item(134, 1).NumberFormat = $dateFormat
$ws.Cells.Item(134, 2).Value = 900200
$ws.Cells.Item(134, 3).Value = 'MSESE'
$ws.Cells.Item(134, 4).Value = 362737.5332642612
$ws.Cells.Item(134, 5).Value = 0
$ws.Cells.Item(134, 6).Value = 0

# Row 135
$ws.Cells.Item(135, 1).Value = 45916
$ws.Cells.Item(135, 1).NumberFormat = $dateFormat
$ws.Cells.Item(135, 2).Value = 900200
$ws.Cells.Item(135, 3).Value = 'SGCIB'
$ws.Cells.Item(135, 4).Value = -164296.2864022301
$ws.Cells.Item(135, 5).Value = 0
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 'cash insuffisant (164 296.29)'

# Row 136
$ws.Cells.Item(136, 1).Value = 45916
$ws.Cells.Item(136, 1).NumberFormat = $dateFormat
$ws.Cells.Item(136, 2).Value = 900200
$ws.Cells.Item(136, 4).Value = 758626386.2456623
$ws.Cells.Item(136, 5).Value = 0
$ws.Cells.Item(136, 6).Value = 0

# Row 137
$ws.Cells.Item(137, 1).Value = 45916
$ws.Cells.Item(137, 1).NumberFormat = $dateFormat
$ws.Cells.Item(137, 2).Value = 981017
$ws.Cells.Item(137, 3).Value = 'BNPP'
$ws.Cells.Item(137, 4).Value = -2061516.132363671
$ws.Cells.Item(137, 5).Value = 0
$ws.Cells.Item(137, 6).Value = 0

# Row 138
$ws.Cells.Item(138, 1).Value = 45916
$ws.Cells.Item(138, 1).NumberFormat = $dateFormat
$ws.Cells.Item(138, 2).Value = 981017
$ws.Cells.Item(138, 3).Value = 'CA'
$ws.Cells.Item(138, 4).Value = -43058.58520903753
$ws.Cells.Item(138, 5).Value = 0
$ws.Cells.Item(138, 6).Value = 0
$ws.Cells.Item(138, 7).Value = 'cash insuffisant (21 328.19)'

# Row 139
$ws.Cells.Item(139, 1).Value = 45916
$ws.Cells.Item(139, 1).NumberFormat = $dateFormat
$ws.Cells.Item(139, 2).Value = 981017
$ws.Cells.Item(139, 3).Value = 'CEP'
$ws.Cells.Item(139, 4).Value = -6289544.212324544
$ws.Cells.Item(139, 5).Value = 0
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 'cash insuffisant (6 289 544.21)'

# Row 140
$ws.Cells.Item(140, 1).Value = 45916
$ws.Cells.Item(140, 1).NumberFormat = $dateFormat
$ws.Cells.Item(140, 2).Value = 981017
$ws.Cells.Item(140, 3).Value = 'DBKAG'
$ws.Cells.Item(140, 4).Value = 507609.1793251648
$ws.Cells.Item(140, 5).Value = 0
$ws.Cells.Item(140, 6).Value = 0

# Row 141
$ws.Cells.Item(141, 1).Value = 45916
$ws.Cells.Item(141, 1).NumberFormat = $dateFormat
$ws.Cells.Item(141, 2).Value = 981017
$ws.Cells.Item(141, 3).Value = 'GIPB'
$ws.Cells.Item(141, 4).Value = -225308.3521440715
$ws.Cells.Item(141, 5).Value = 0
$ws.Cells.Item(141, 6).Value = 0
$ws.Cells.Item(141, 7).Value = 'cash insuffisant (225 308.35)'

# Row 142
$ws.Cells.Item(142, 1).Value = 45916
$ws.Cells.Item(142, 1).NumberFormat = $dateFormat
$ws.Cells.Item(142, 2).Value = 981017
$ws.Cells.Item(142, 3).Value = 'JPMSE'
$ws.Cells.Item(142, 4).Value = -1105152.639062842
$ws.Cells.Item(142, 5).Value = 0
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 'cash insuffisant (1 105 152.64)'

# Row 143
$ws.Cells.Item(143, 1).Value = 45916
$ws.Cells.Item(143, 1).NumberFormat = $dateFormat
$ws.Cells.Item(143, 2).Value = 981017
$ws.Cells.Item(143, 3).Value = 'MSESE'
$ws.Cells.Item(143, 4).Value = -26860711.45644478
$ws.Cells.Item(143, 5).Value = 0
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 'cash insuffisant (26 860 711.46)'

# Row 144
$ws.Cells.Item(144, 1).Value = 45916
$ws.Cells.Item(144, 1).NumberFormat = $dateFormat
$ws.Cells.Item(144, 2).Value = 981017
$ws.Cells.Item(144, 3).Value = 'NATIXIS'
$ws.Cells.Item(144, 4).Value = -41914016.5082352
$ws.Cells.Item(144, 5).Value = 0
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 'cash insuffisant (41 914 016.51)'

# Row 145
$ws.Cells.Item(145, 1).Value = 45916
$ws.Cells.Item(145, 1).NumberFormat = $dateFormat
$ws.Cells.Item(145, 2).Value = 981017
$ws.Cells.Item(145, 4).Value = 417571338.1490759
$ws.Cells.Item(145, 5).Value = 0
$ws.Cells.Item(145, 6).Value = 0

# Row 146
$ws.Cells.Item(146, 1).Value = 45916
$ws.Cells.Item(146, 1).NumberFormat = $dateFormat
$ws.Cells.Item(146, 2).Value = 900200
$ws.Cells.Item(146, 3).Value = 'BNP PARIBAS'
$ws.Cells.Item(146, 4).Value = -44748309.06889844
$ws.Cells.Item(146, 5).Value = 0
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 'cash insuffisant (22 757 791.25)'

# Row 147
$ws.Cells.Item(147, 1).Value = 45916
$ws.Cells.Item(147, 1).NumberFormat = $dateFormat
$ws.Cells.Item(147, 2).Value = 900200
$ws.Cells.Item(147, 3).Value = 'CITIBANK EUROPE PUBLIC LIMITED COMPANY'
$ws.Cells.Item(147, 4).Value = 110691.9541025327
$ws.Cells.Item(147, 5).Value = 0
$ws.Cells.Item(147, 6).Value = 0

# Row 148
$ws.Cells.Item(148, 1).Value = 45916
$ws.Cells.Item(148, 1).NumberFormat = $dateFormat
$ws.Cells.Item(148, 2).Value = 900200
$ws.Cells.Item(148, 3).Value = 'CREDIT AGRICOLE CORPORATE AND INVESTMENT BANK'
$ws.Cells.Item(148, 4).Value = -16145545.642625
$ws.Cells.Item(148, 5).Value = 0
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 'cash insuffisant (16 145 545.64)'

# Row 149
$ws.Cells.Item(149, 1).Value = 45916
$ws.Cells.Item(149, 1).NumberFormat = $dateFormat
$ws.Cells.Item(149, 2).Value = 900200
$ws.Cells.Item(149, 3).Value = 'GOLDMAN SACHS BANK EUROPE SE'
$ws.Cells.Item(149, 4).Value = -1552465.916127797
$ws.Cells.Item(149, 5).Value = 0
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 'cash insuffisant (1 552 465.92)'

# Row 150
$ws.Cells.Item(150, 1).Value = 45916
$ws.Cells.Item(150, 1).NumberFormat = $dateFormat
$ws.Cells.Item(150, 2).Value = 900200
$ws.Cells.Item(150, 3).Value = 'GOLDMAN SACHS INTERNATIONAL  PARIS BRANCH'
$ws.Cells.Item(150, 4).Value = -262844.824426329
$ws.Cells.Item(150, 5).Value = 0
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 'cash insuffisant (262 844.82)'

# Row 151
$ws.Cells.Item(151, 1).Value = 45916
$ws.Cells.Item(151, 1).NumberFormat = $dateFormat
$ws.Cells.Item(151, 2).Value = 900200
$ws.Cells.Item(151, 3).Value = 'JP MORGAN SE'
$ws.Cells.Item(151, 4).Value = -6725967.70460027
$ws.Cells.Item(151, 5).Value = 0
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 'cash insuffisant (6 725 967.70)'

# Row 152
$ws.Cells.Item(152, 1).Value = 45916
$ws.Cells.Item(152, 1).NumberFormat = $dateFormat
$ws.Cells.Item(152, 2).Value = 900200
$ws.Cells.Item(152, 3).Value = 'MORGAN STANLEY EUROPE SE'
$ws.Cells.Item(152, 4).Value = 362737.5332642612
$ws.Cells.Item(152, 5).Value = 0
$ws.Cells.Item(152, 6).Value = 0

# Row 153
$ws.Cells.Item(153, 1).Value = 45916
$ws.Cells.Item(153, 1).NumberFormat = $dateFormat
$ws.Cells.Item(153, 2).Value = 900200
$ws.Cells.Item(153, 3).Value = 'SOCIETE GENERALE'
$ws.Cells.Item(153, 4).Value = -164296.2864022301
$ws.Cells.Item(153, 5).Value = 0
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 'cash insuffisant (164 296.29)'

# Row 154
$ws.Cells.Item(154, 1).Value = 45916
$ws.Cells.Item(154, 1).NumberFormat = $dateFormat
$ws.Cells.Item(154, 2).Value = 900200
$ws.Cells.Item(154, 4).Value = 758626386.2456623
$ws.Cells.Item(154, 5).Value = 0
$ws.Cells.Item(154, 6).Value = 0

# Row 155
$ws.Cells.Item(155, 1).Value = 45916
$ws.Cells.Item(155, 1).NumberFormat = $dateFormat
$ws.Cells.Item(155, 2).Value = 981017
$ws.Cells.Item(155, 3).Value = 'BNP PARIBAS'
$ws.Cells.Item(155, 4).Value = -2061516.132363671
$ws.Cells.Item(155, 5).Value = 0
$ws.Cells.Item(155, 6).Value = 0

# Row 156
$ws.Cells.Item(156, 1).Value = 45916
$ws.Cells.Item(156, 1).NumberFormat = $dateFormat
$ws.Cells.Item(156, 2).Value = 981017
$ws.Cells.Item(156, 3).Value = 'CITIBANK EUROPE PUBLIC LIMITED COMPANY'
$ws.Cells.Item(156, 4).Value = -6289544.212324544
$ws.Cells.Item(156, 5).Value = 0
$ws.Cells.Item(156, 6).Value = 0
$ws.Cells.Item(156, 7).Value = 'cash insuffisant (6 267 813.82)'

# Row 157
$ws.Cells.Item(157, 1).Value = 45916
$ws.Cells.Item(157, 1).NumberFormat = $dateFormat
$ws.Cells.Item(157, 2).Value = 981017
$ws.Cells.Item(157, 3).Value = 'CREDIT AGRICOLE CORPORATE AND INVESTMENT BANK'
$ws.Cells.Item(157, 4).Value = -43058.58520903753
$ws.Cells.Item(157, 5).Value = 0
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 'cash insuffisant (43 058.59)'

# Row 158
$ws.Cells.Item(158, 1).Value = 45916
$ws.Cells.Item(158, 1).NumberFormat = $dateFormat
$ws.Cells.Item(158, 2).Value = 981017
$ws.Cells.Item(158, 3).Value = 'DEUTSCHE BANK AKTIENGESELLSCHAFT'
$ws.Cells.Item(158, 4).Value = 507609.1793251648
$ws.Cells.Item(158, 5).Value = 0
$ws.Cells.Item(158, 6).Value = 0

# Row 159
$ws.Cells.Item(159, 1).Value = 45916
$ws.Cells.Item(159, 1).NumberFormat = $dateFormat
$ws.Cells.Item(159, 2).Value = 981017
$ws.Cells.Item(159, 3).Value = 'GOLDMAN SACHS INTERNATIONAL  PARIS BRANCH'
$ws.Cells.Item(159, 4).Value = -225308.3521440715
$ws.Cells.Item(159, 5).Value = 0
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 'cash insuffisant (225 308.35)'

# Row 160
$ws.Cells.Item(160, 1).Value = 45916
$ws.Cells.Item(160, 1).NumberFormat = $dateFormat
$ws.Cells.Item(160, 2).Value = 981017
$ws.Cells.Item(160, 3).Value = 'JP MORGAN SE'
$ws.Cells.Item(160, 4).Value = -1105152.639062842
$ws.Cells.Item(160, 5).Value = 0
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 'cash insuffisant (1 105 152.64)'

# Row 161
$ws.Cells.Item(161, 1).Value = 45916
$ws.Cells.Item(161, 1).NumberFormat = $dateFormat
$ws.Cells.Item(161, 2).Value = 981017
$ws.Cells.Item(161, 3).Value = 'MORGAN STANLEY EUROPE SE'
$ws.Cells.Item(161, 4).Value = -26860711.45644478
$ws.Cells.Item(161, 5).Value = 0
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 'cash insuffisant (26 860 711.46)'

# Row 162
$ws.Cells.Item(162, 1).Value = 45916
$ws.Cells.Item(162, 1).NumberFormat = $dateFormat
$ws.Cells.Item(162, 2).Value = 981017
$ws.Cells.Item(162, 3).Value = 'NATIXIS'
$ws.Cells.Item(162, 4).Value = -41914016.5082352
$ws.Cells.Item(162, 5).Value = 0
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 'cash insuffisant (41 914 016.51)'

# Row 163
$ws.Cells.Item(163, 1).Value = 45916
$ws.Cells.Item(163, 1).NumberFormat = $dateFormat
$ws.Cells.Item(163, 2).Value = 981017
$ws.Cells.Item(163, 4).Value = 417571338.1490759
$ws.Cells.Item(163, 5).Value = 0
$ws.Cells.Item(163, 6).Value = 0

# Row 164
$ws.Cells.Item(164, 1).Value = 45916
$ws.Cells.Item(164, 1).NumberFormat = $dateFormat
$ws.Cells.Item(164, 2).Value = 900200
$ws.Cells.Item(164, 3).Value = 'BNP PARIBAS'
$ws.Cells.Item(164, 4).Value = -44748309.06889844
$ws.Cells.Item(164, 5).Value = 0
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 'cash insuffisant (22 757 791.25)'

# Row 165
$ws.Cells.Item(165, 1).Value = 45916
$ws.Cells.Item(165, 1).NumberFormat = $dateFormat
$ws.Cells.Item(165, 2).Value = 900200
$ws.Cells.Item(165, 3).Value = 'CITIBANK EUROPE PUBLIC LIMITED COMPANY'
$ws.Cells.Item(165, 4).Value = 110691.9541025327
$ws.Cells.Item(165, 5).Value = 0
$ws.Cells.Item(165, 6).Value = 0

# Row 166
$ws.Cells.Item(166, 1).Value = 45916
$ws.Cells.Item(166, 1).NumberFormat = $dateFormat
$ws.Cells.Item(166, 2).Value = 900200
$ws.Cells.Item(166, 3).Value = 'CREDIT AGRICOLE CORPORATE AND INVESTMENT BANK'
$ws.Cells.Item(166, 4).Value = -16145545.642625
$ws.Cells.Item(166, 5).Value = 0
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 'cash insuffisant (16 145 545.64)'

# Row 167
$ws.Cells.Item(167, 1).Value = 45916
$ws.Cells.Item(167, 1).NumberFormat = $dateFormat
$ws.Cells.Item(167, 2).Value = 900200
$ws.Cells.Item(167, 3).Value = 'GOLDMAN SACHS BANK EUROPE SE'
$ws.Cells.Item(167, 4).Value = -1552465.916127797
$ws.Cells.Item(167, 5).Value = 0
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 'cash insuffisant (1 552 465.92)'

# Row 168
$ws.Cells.Item(168, 1).Value = 45916
$ws.Cells.Item(168, 1).NumberFormat = $dateFormat
$ws.Cells.Item(168, 2).Value = 900200
$ws.Cells.Item(168, 3).Value = 'GOLDMAN SACHS INTERNATIONAL  PARIS BRANCH'
$ws.Cells.Item(168, 4).Value = -262844.824426329
$ws.Cells.Item(168, 5).Value = 0
$ws.Cells.Item(168, 6).Value = 0
$ws.Cells.Item(168, 7).Value = 'cash insuffisant (262 844.82)'

# Row 169
$ws.Cells.Item(169, 1).Value = 45916
$ws.Cells.Item(169, 1).NumberFormat = $dateFormat
$ws.Cells.Item(169, 2).Value = 900200
$ws.Cells.Item(169, 3).Value = 'JP MORGAN SE'
$ws.Cells.Item(169, 4).Value = -6725967.70460027
$ws.Cells.Item(169, 5).Value = 0
$ws.Cells.Item(169, 6).Value = 0
$ws.Cells.Item(169, 7).Value = 'cash insuffisant (6 725 967.70)'

# Row 170
$ws.Cells.Item(170, 1).Value = 45916
$ws.Cells.Item(170, 1).NumberFormat = $dateFormat
$ws.Cells.Item(170, 2).Value = 900200
$ws.Cells.Item(170, 3).Value = 'MORGAN STANLEY EUROPE SE'
$ws.Cells.Item(170, 4).Value = 362737.5332642612
$ws.Cells.Item(170, 5).Value = 0
$ws.Cells.Item(170, 6).Value = 0

# Row 171
$ws.Cells.Item(171, 1).Value = 45916
$ws.Cells.Item(171, 1).NumberFormat = $dateFormat
$ws.Cells.Item(171, 2).Value = 900200
$ws.Cells.Item(171, 3).Value = 'SOCIETE GENERALE'
$ws.Cells.Item(171, 4).Value = -164296.2864022301
$ws.Cells.Item(171, 5).Value = 0
$ws.Cells.Item(171, 6).Value = 0
$ws.Cells.Item(171, 7).Value = 'cash insuffisant (164 296.29)'

# Row 172
$ws.Cells.Item(172, 1).Value = 45916
$ws.Cells.Item(172, 1).NumberFormat = $dateFormat
$ws.Cells.Item(172, 2).Value = 900200
$ws.Cells.Item(172, 4).Value = 758626386.2456623
$ws.Cells.Item(172, 5).Value = 0
$ws.Cells.Item(172, 6).Value = 0

# Row 173
$ws.Cells.Item(173, 1).Value = 45916
$ws.Cells.Item(173, 1).NumberFormat = $dateFormat
$ws.Cells.Item(173, 2).Value = 981017
$ws.Cells.Item(173, 3).Value = 'BNP PARIBAS'
$ws.Cells.Item(173, 4).Value = -2061516.132363671
$ws.Cells.Item(173, 5).Value = 0
$ws.Cells.Item(173, 6).Value = 0

# Row 174
$ws.Cells.Item(174, 1).Value = 45916
$ws.Cells.Item(174, 1).NumberFormat = $dateFormat
$ws.Cells.Item(174, 2).Value = 981017
$ws.Cells.Item(174, 3).Value = 'CITIBANK EUROPE PUBLIC LIMITED COMPANY'
$ws.Cells.Item(174, 4).Value = -6289544.212324544
$ws.Cells.Item(174, 5).Value = 0
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 'cash insuffisant (6 267 813.82)'

# Row 175
$ws.Cells.Item(175, 1).Value = 45916
$ws.Cells.Item(175, 1).NumberFormat = $dateFormat
$ws.Cells.Item(175, 2).Value = 981017
$ws.Cells.Item(175, 3).Value = 'CREDIT AGRICOLE CORPORATE AND INVESTMENT BANK'
$ws.Cells.Item(175, 4).Value = -43058.58520903753
$ws.Cells.Item(175, 5).Value = 0
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 'cash insuffisant (43 058.59)'

# Row 176
$ws.Cells.Item(176, 1).Value = 45916
$ws.Cells.Item(176, 1).NumberFormat = $dateFormat
$ws.Cells.Item(176, 2).Value = 981017
$ws.Cells.Item(176, 3).Value = 'DEUTSCHE BANK AKTIENGESELLSCHAFT'
$ws.Cells.Item(176, 4).Value = 507609.1793251648
$ws.Cells.Item(176, 5).Value = 0
$ws.Cells.Item(176, 6).Value = 0

# Row 177
$ws.Cells.Item(177, 1).Value = 45916
$ws.Cells.Item(177, 1).NumberFormat = $dateFormat
$ws.Cells.Item(177, 2).Value = 981017
$ws.Cells.Item(177, 3).Value = 'GOLDMAN SACHS INTERNATIONAL  PARIS BRANCH'
$ws.Cells.Item(177, 4).Value = -225308.3521440715
$ws.Cells.Item(177, 5).Value = 0
$ws.Cells.Item(177, 6).Value = 0
$ws.Cells.Item(177, 7).Value = 'cash insuffisant (225 308.35)'

# Row 178
$ws.Cells.Item(178, 1).Value = 45916
$ws.Cells.Item(178, 1).NumberFormat = $dateFormat
$ws.Cells.Item(178, 2).Value = 981017
$ws.Cells.Item(178, 3).Value = 'JP MORGAN SE'
$ws.Cells.Item(178, 4).Value = -1105152.639062842
$ws.Cells.Item(178, 5).Value = 0
$ws.Cells.Item(178, 6).Value = 0
$ws.Cells.Item(178, 7).Value = 'cash insuffisant (1 105 152.64)'

# Row 179
$ws.Cells.Item(179, 1).Value = 45916
$ws.Cells.Item(179, 1).NumberFormat = $dateFormat
$ws.Cells.Item(179, 2).Value = 981017
$ws.Cells.Item(179, 3).Value = 'MORGAN STANLEY EUROPE SE'
$ws.Cells.Item(179, 4).Value = -26860711.45644478
$ws.Cells.Item(179, 5).Value = 0
$ws.Cells.Item(179, 6).Value = 0
$ws.Cells.Item(179, 7).Value = 'cash insuffisant (26 860 711.46)'

# Row 180
$ws.Cells.Item(180, 1).Value = 45916
$ws.Cells.Item(180, 1).NumberFormat = $dateFormat
$ws.Cells.Item(180, 2).Value = 981017
$ws.Cells.Item(180, 3).Value = 'NATIXIS'
$ws.Cells.Item(180, 4).Value = -41914016.5082352
$ws.Cells.Item(180, 5).Value = 0
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 'cash insuffisant (41 914 016.51)'

# Row 181
$ws.Cells.Item(181, 1).Value = 45916
$ws.Cells.Item(181, 1).NumberFormat = $dateFormat
$ws.Cells.Item(181, 2).Value = 981017
$ws.Cells.Item(181, 4).Value = 417571338.1490759
$ws.Cells.Item(181, 5).Value = 0
$ws.Cells.Item(181, 6).Value = 0

